$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Professional summary: neutralize "all Black and Asian-American voters"
#    -> "50M voters" (plain text, no formatting change)
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: split "all Black and Asian-American voters" so
#    "50M" becomes its own bold run (matching the styling used for the other
#    stat callouts in that same bullet).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$find2 = $rng.Find
$find2.ClearFormatting()
$find2.Execute("all Black and Asian-American", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "50M"
$rng.Font.Bold = $true
$rng.Font.Color = 5258796   # RGB 0x2C3E50 encoded as BGR for Word's Font.Color

# ---------------------------------------------------------------------------
# 3) Reorder work-experience entries: move "Research Director - PCCC" and
#    "Data Products Manager - Helm/Murmuration" from the end of the
#    experience list (after the Praxis Project entry) up to sit right after
#    "Partner - Siege Analytics" / around "Software Engineer - Mautinoa".
#    Final order: ... Siege Analytics -> Data Products Manager ->
#    Software Engineer - Mautinoa -> Research Director - PCCC ->
#    Software Engineer - Salsa Labs -> ... -> Interim Technology Manager.
# ---------------------------------------------------------------------------

# Locate the two heading paragraphs that bound the block to move.
$researchDirIdx = 0
$praxisBulletIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Research Director - PCCC*") { $researchDirIdx = $i }
    if ($t -like "*Implemented CRM systems for stakeholder engagement*") { $praxisBulletIdx = $i }
}
# The block spans from the "Research Director" heading through the end of
# the "Data Products Manager" bullets, i.e. the 10 paragraphs that precede
# "KEY PROJECTS".
$keyProjectsIdx = 0
for ($i = $researchDirIdx; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "KEY PROJECTS*") { $keyProjectsIdx = $i; break }
}
$blockStart = $d.Paragraphs($researchDirIdx).Range.Start
$blockEnd = $d.Paragraphs($keyProjectsIdx - 1).Range.End
$d.Range($blockStart, $blockEnd).Delete()

# Find "Software Engineer - Mautinoa Technologies" and "Geospatial analysis
# on populations..." (last bullet of that entry) after the deletion above.
$mautinoaIdx = 0
$geospatialBulletIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Software Engineer - Mautinoa Technologies*") { $mautinoaIdx = $i }
    if ($t -like "*Geospatial analysis on populations and boundaries*") { $geospatialBulletIdx = $i }
}

# --- Insert "Data Products Manager" entry right before Mautinoa ---------
$anchor = $d.Paragraphs($mautinoaIdx)
$anchor.Range.InsertParagraphBefore() | Out-Null
$heading = $d.Paragraphs($mautinoaIdx)
$heading.Range.Text = "Data Products Manager - Helm/Murmuration (Austin, TX) | 2021 - 2023"

$anchor = $d.Paragraphs($mautinoaIdx)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($mautinoaIdx + 1)
$p.Style = "Normal"
$p.Range.Text = "Democratic Electoral Technology"

$anchor = $d.Paragraphs($mautinoaIdx + 1)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($mautinoaIdx + 2)
$p.Range.Text = "• Led design and implementation of enterprise-scale multi-tenant data warehouse for geo-referenced demographic, econometric, and electoral data"

$anchor = $d.Paragraphs($mautinoaIdx + 2)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($mautinoaIdx + 3)
$p.Range.Text = "• Managed engineering team of 11 professionals while setting technical direction for data architecture"

$anchor = $d.Paragraphs($mautinoaIdx + 3)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($mautinoaIdx + 4)
$p.Range.Text = "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by "
$etlRange = $p.Range
$etlRange.Collapse(0)
$etlRange.InsertAfter("57%")
$etlRange.Font.Bold = $true
$etlRange.Font.Color = 5258796

# "Mautinoa" entry (and everything after) has shifted down by 5 paragraphs.
$mautinoaIdx = $mautinoaIdx + 5
$geospatialBulletIdx = $geospatialBulletIdx + 5

# --- Insert "Research Director - PCCC" entry right after the Mautinoa
#     entry's last bullet ("Geospatial analysis on populations...") --------
$anchor = $d.Paragraphs($geospatialBulletIdx)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($geospatialBulletIdx + 1)
$p.Style = "Heading 3"
$p.Range.Text = "Research Director - PCCC (Washington, DC) | August 2011 - August 2012"

$anchor = $d.Paragraphs($geospatialBulletIdx + 1)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($geospatialBulletIdx + 2)
$p.Style = "Normal"
$p.Range.Text = "Political Research & Data Analysis (FLEEM System)"

$anchor = $d.Paragraphs($geospatialBulletIdx + 2)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($geospatialBulletIdx + 3)
$p.Range.Text = "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys"

$anchor = $d.Paragraphs($geospatialBulletIdx + 3)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($geospatialBulletIdx + 4)
$p.Range.Text = "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"

$anchor = $d.Paragraphs($geospatialBulletIdx + 4)
$anchor.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs($geospatialBulletIdx + 5)
$p.Range.Text = "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
